$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-12-13"

# Update the header cell text for the "2022 (through ...)" column
$ws.Range("I1").Value = "2022 (through 12-13)"

# Update the December figure and the Total figure for that column
$ws.Range("I13").Value = 57
$ws.Range("I14").Value = 1573
